# "added 4wk low sales check"
# Updates the forecast figures (MyForecast, Inventory Coverage, Seasonality
# Index and the derived risk/urgency labels) on the "Forecast Comparison"
# sheet to reflect the new 4-week low-sales check, and refreshes the
# dependent roll-up figures on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------
# MyForecast (D), Inventory Coverage (H), Seasonality Index (L) per week,
# plus the Stockout Risk (I) / Reorder Urgency (J) labels that flip once
# inventory coverage runs out.

$rows = @(
    @{ Row = 2;  D = 57; H = 9.25;               L = 1.02 },
    @{ Row = 3;  D = 57; H = 8.25;               L = 1.07 },
    @{ Row = 4;  D = 56; H = 7.38;               L = 0.98 },
    @{ Row = 5;  D = 54; H = 6.61;               L = 1.04 },
    @{ Row = 6;  D = 51; H = 5.94;               L = 1.1  },
    @{ Row = 7;  D = 49; H = 5.14;               L = 0.96 },
    @{ Row = 8;  D = 47; H = 4.32;               L = 0.82 },
    @{ Row = 9;        H = 3.47;               L = 0.9  },
    @{ Row = 10;       H = 2.64;               L = 1.12 },
    @{ Row = 11;       H = 1.64;               L = 1.03 },
    @{ Row = 12; D = 39; H = 0.6899999999999999; L = 1.16; J = "Urgent" },
    @{ Row = 13; D = 37; H = 0;                 L = 1.07; I = "High"; J = "Urgent" },
    @{ Row = 14; D = 35; H = 0;                 L = 1.12 },
    @{ Row = 15; D = 36;                        L = 0.9399999999999999 },
    @{ Row = 16; D = 33;                        L = 1    },
    @{ Row = 17;                                L = 0.87 }
)

foreach ($r in $rows) {
    $row = $r.Row
    if ($r.ContainsKey("D")) { $wsForecast.Cells.Item($row, 4).Value  = $r.D }
    if ($r.ContainsKey("H")) { $wsForecast.Cells.Item($row, 8).Value  = $r.H }
    if ($r.ContainsKey("I")) { $wsForecast.Cells.Item($row, 9).Value  = $r.I }
    if ($r.ContainsKey("J")) { $wsForecast.Cells.Item($row, 10).Value = $r.J }
    if ($r.ContainsKey("L")) { $wsForecast.Cells.Item($row, 12).Value = $r.L }
}

# --- Summary sheet ---------------------------------------------------------
# Roll-up totals recompute after the forecast refresh above. These are
# stored as text on the sheet, so re-apply the text formatting (a leading
# apostrophe, same as typing the value into Excel) to keep them as text
# rather than letting them coerce to numbers, then drop the resulting
# "number stored as text" formatting so the cell format is left untouched.

$summaryUpdates = @(
    @{ Cell = "B9";  Value = "713" },  # Total Forecast (16 Weeks)
    @{ Cell = "B10"; Value = "416" },  # Total Forecast (8 Weeks)
    @{ Cell = "B11"; Value = "224" },  # Total Forecast (4 Weeks)
    @{ Cell = "B12"; Value = "57"  }   # Max Forecast
)

foreach ($u in $summaryUpdates) {
    $cell = $wsSummary.Range($u.Cell)
    $cell.Value = "'" + $u.Value
    $cell.ClearFormats()
}
